# Estadisticos Segundo Parcial 26 Mayo
#
# "Rescatables" sheet: register a newly-identified resit student and refresh
# the failure counts, then keep the list sorted the way the report expects
# (most-failed first, grouped by class).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Make room for the new student right below the first record.
$ws.Rows.Item(3).Insert()

$ws.Cells.Item(3, 1).Value = 24330051920187
$ws.Cells.Item(3, 2).Value = "OSORIO"
$ws.Cells.Item(3, 3).Value = "HERNANDEZ"
$ws.Cells.Item(3, 4).Value = "AYLIN ABIGAIL"
$ws.Cells.Item(3, 5).Value = "GESTIONA DOCUMENTACIÓN DEL ÁREA DE RECURSOS HUMANOS"
$ws.Cells.Item(3, 6).Value = "2ARHV"
$ws.Cells.Item(3, 7).Value = 4

# Refreshed Reprobadas (failed-subject) tallies for three existing students,
# now sitting one row lower because of the insert above.
$ws.Cells.Item(5, 7).Value = 3   # OREA MARTINEZ, JOSE MANUEL: 4 -> 3
$ws.Cells.Item(6, 7).Value = 3   # PALMA RANGEL, ROBERTO: 4 -> 3
$ws.Cells.Item(11, 7).Value = 2  # REYES TLAXCALTECA, GAEL ARMANDO: 3 -> 2

# Re-sort the whole list: most Reprobadas first, then by Grupo.
$sortRange = $ws.Range("A2:G16")
$keyReprobadas = $ws.Range("G2:G16")
$keyGrupo = $ws.Range("F2:F16")
$sortRange.Sort($keyReprobadas, 2, $keyGrupo, $null, 1)
